# Sync attendance_reports: fix "Recorded By" (column G) ordering so that
# "System" is listed last instead of first in the comma-separated list.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dim = $ws.UsedRange
$lastRow = $dim.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2
    if ($null -ne $val -and $val.ToString().StartsWith("System, ")) {
        $parts = $val.ToString().Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $first = $parts[0]
        $rest = $parts[1..($parts.Length - 1)]
        $newParts = $rest + ,$first
        $newVal = [string]::Join(", ", $newParts)
        $cell.Value = $newVal
    }
}
